$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 46073

# Row 3
$ws.Range("C3").Value = 46073

# Row 4
$ws.Range("A4").Value = "A 39958-2024"
$ws.Range("B4").Value = 45553
$ws.Range("C4").Value = 46073
$ws.Range("G4").Value = 3.4

# Row 5
$ws.Range("A5").Value = "A 39924-2025"
$ws.Range("C5").Value = 46073
$ws.Range("G5").Value = 1.1

# Row 6
$ws.Range("A6").Value = "A 39928-2025"
$ws.Range("C6").Value = 46073
$ws.Range("G6").Value = 0.9

# Row 7
$ws.Range("A7").Value = "A 40001-2025"
$ws.Range("B7").Value = 45891
$ws.Range("C7").Value = 46073
$ws.Range("G7").Value = 0.6

# Row 8
$ws.Range("A8").Value = "A 2229-2023"
$ws.Range("B8").Value = 44939
$ws.Range("C8").Value = 46073
$ws.Range("F8").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G8").Value = 4.3

# Row 9
$ws.Range("C9").Value = 46073

# Row 10
$ws.Range("A10").Value = "A 32256-2025"
$ws.Range("B10").Value = 45835.6353125
$ws.Range("C10").Value = 46073
$ws.Range("F10").ClearContents()
$ws.Range("G10").Value = 5.4

# Row 11
$ws.Range("A11").Value = "A 39876-2024"
$ws.Range("B11").Value = 45553
$ws.Range("C11").Value = 46073
$ws.Range("G11").Value = 0.3

# Row 12
$ws.Range("A12").Value = "A 7694-2023"
$ws.Range("B12").Value = 44967
$ws.Range("C12").Value = 46073
$ws.Range("G12").Value = 2.2

# Row 13
$ws.Range("A13").Value = "A 5968-2023"
$ws.Range("B13").Value = 44959
$ws.Range("C13").Value = 46073
$ws.Range("G13").Value = 1.5

# Row 14
$ws.Range("A14").Value = "A 6004-2026"
$ws.Range("B14").Value = 46050
$ws.Range("C14").Value = 46073
$ws.Range("G14").Value = 2.7

# Row 15
$ws.Range("A15").Value = "A 28815-2024"
$ws.Range("B15").Value = 45478
$ws.Range("C15").Value = 46073
$ws.Range("G15").Value = 2.8

# Row 16
$ws.Range("A16").Value = "A 34926-2022"
$ws.Range("B16").Value = 44796
$ws.Range("C16").Value = 46073
$ws.Range("G16").Value = 1.3

# Row 17
$ws.Range("A17").Value = "A 2727-2024"
$ws.Range("B17").Value = 45314
$ws.Range("C17").Value = 46073
$ws.Range("G17").Value = 3.8

# Row 18
$ws.Range("A18").Value = "A 5528-2023"
$ws.Range("B18").Value = 44957
$ws.Range("C18").Value = 46073
$ws.Range("G18").Value = 1.2

